# Update course schedule for new term (Winter 2017)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "cm"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "link_it"
$ws.Range("D1").Value = "topic"

$ws.Range("A2").Value = "cm001"
$ws.Range("B2").Value = 42739
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "Introduction to computing for the social sciences"

$ws.Range("A3").Value = "cm002"
$ws.Range("B3").Formula = "=B2+5"
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = "Visualizations and the grammar of graphics"

$ws.Range("A4").Value = "cm003"
$ws.Range("B4").Formula = "=B2+7"
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = "Data transformation and exploratory data analysis"

$ws.Range("A5").Value = "cm000"
$ws.Range("B5").Formula = "=B4+5"
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = "No class (MLK Jr. Day)"

$ws.Range("A6").Value = "cm004"
$ws.Range("B6").Formula = "=B4+7"
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = "Data wrangling"

$ws.Range("A7").Value = "cm005"
$ws.Range("B7").Formula = "=B6+5"
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = "Pipes and functions in R"

$ws.Range("A8").Value = "cm006"
$ws.Range("B8").Formula = "=B6+7"
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = "Vectors and iteration"

$ws.Range("A9").Value = "cm007"
$ws.Range("B9").Formula = "=B8+5"
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = "Debugging errors"

$ws.Range("A10").Value = "cm008"
$ws.Range("B10").Formula = "=B8+7"
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = "Reproducibility in research"

$ws.Range("A11").Value = "cm009"
$ws.Range("B11").Formula = "=B10+5"
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = "Statistical learning: regression"

$ws.Range("A12").Value = "cm010"
$ws.Range("B12").Formula = "=B10+7"
$ws.Range("C12").Value = $true
$ws.Range("D12").Value = "Statistical learning: classification"

$ws.Range("A13").Value = "cm011"
$ws.Range("B13").Formula = "=B12+5"
$ws.Range("C13").Value = $true
$ws.Range("D13").Value = "Statistical learning: cross-validation"

$ws.Range("A14").Value = "cm012"
$ws.Range("B14").Formula = "=B12+7"
$ws.Range("C14").Value = $true
$ws.Range("D14").Value = "Distributed computing"

$ws.Range("A15").Value = "cm013"
$ws.Range("B15").Formula = "=B14+5"
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = "Getting data from the web: API access"

$ws.Range("A16").Value = "cm014"
$ws.Range("B16").Formula = "=B14+7"
$ws.Range("C16").Value = $true
$ws.Range("D16").Value = "Getting data from the web: scraping"

$ws.Range("A17").Value = "cm015"
$ws.Range("B17").Formula = "=B16+5"
$ws.Range("C17").Value = $true
$ws.Range("D17").Value = "Text analysis: fundamentals and sentiment analysis"

$ws.Range("A18").Value = "cm016"
$ws.Range("B18").Formula = "=B16+7"
$ws.Range("C18").Value = $true
$ws.Range("D18").Value = "Text analysis: topic modeling"

$ws.Range("A19").Value = "cm017"
$ws.Range("B19").Formula = "=B18+5"
$ws.Range("C19").Value = $true
$ws.Range("D19").Value = "Building Shiny apps"

$ws.Range("A20").Value = "cm018"
$ws.Range("B20").Formula = "=B18+7"
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = "Building Shiny apps (part II)"

$ws.Range("A21").Value = "lab01"
$ws.Range("B21").Value = 42739
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = "Software setup"

$ws.Range("A22").Value = "lab02"
$ws.Range("B22").Value = 42744
$ws.Range("C22").Value = $false
$ws.Range("D22").ClearContents()

$ws.Range("A23").Value = "lab03"
$ws.Range("B23").Formula = "=B22+7"
$ws.Range("C23").Value = $false
$ws.Range("D23").ClearContents()

$ws.Range("A24").Value = "lab04"
$ws.Range("B24").Formula = "=B23+7"
$ws.Range("C24").Value = $false
$ws.Range("D24").ClearContents()

$ws.Range("A25").Value = "lab05"
$ws.Range("B25").Formula = "=B24+7"
$ws.Range("C25").Value = $false
$ws.Range("D25").ClearContents()

$ws.Range("A26").Value = "lab06"
$ws.Range("B26").Formula = "=B25+7"
$ws.Range("C26").Value = $false
$ws.Range("D26").ClearContents()

$ws.Range("A27").Value = "lab07"
$ws.Range("B27").Formula = "=B26+7"
$ws.Range("C27").Value = $false
$ws.Range("D27").ClearContents()

$ws.Range("A28").Value = "lab08"
$ws.Range("B28").Formula = "=B27+7"
$ws.Range("C28").Value = $false
$ws.Range("D28").ClearContents()

$ws.Range("A29").Value = "lab09"
$ws.Range("B29").Formula = "=B28+7"
$ws.Range("C29").Value = $false
$ws.Range("D29").ClearContents()

$ws.Range("A30").Value = "lab10"
$ws.Range("B30").Formula = "=B29+7"
$ws.Range("C30").Value = $false
$ws.Range("D30").ClearContents()
# Remove the now-unused 31st row (term now has 30 rows)
$ws.Rows.Item(31).Delete()

# Update the remembered selection to match the new sheet state
[void]$ws.Range("A6").Select()
